$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Table 2.3 (sheet1) edits
# ---------------------------------------------------------------------------

# A2 was the text "boreal" -- it becomes literal 0.
$ws1.Cells.Item(2,1).Value = 0

# G12 precision tweak
$ws1.Cells.Item(12,7).Value = 36.8333333333333

# New I12 formula: average of G12 and G13
$ws1.Cells.Item(12,9).Formula = "=(G12+G13)/2"

# B17 used to be a formula; now it is a single literal space character,
# which causes the dependent E17/F17 ratio formulas to error out (#VALUE!)
$ws1.Cells.Item(17,2).Value = " "

# G17 precision tweak
$ws1.Cells.Item(17,7).Value = 60.190476190476197

# New I17 formula: average of G15, G16, G14
$ws1.Cells.Item(17,9).Formula = "=(G15+G16+G14)/3"

# Move the selection on Table 2.3 (it will lose tab focus once Table 5.5
# is created/activated below)
$ws1.Range("G13").Select()

# ---------------------------------------------------------------------------
# Add the new "Table 5.5" worksheet right after "Table 2.3"
# ---------------------------------------------------------------------------

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Table 5.5"

# Row 1
$ws2.Cells.Item(1,1).Value = "LU Fctors"

# Row 2 - FLU / Long Term Cultivated / Cool Temperate/Boreal / dry
$ws2.Cells.Item(2,1).Value = "FLU"
$ws2.Cells.Item(2,2).Value = "Long Term Cultivated"
$ws2.Cells.Item(2,3).Value = "Cool Temperate/Boreal"
$ws2.Cells.Item(2,4).Value = "-"
$ws2.Cells.Item(2,5).Value = 0.77

# Row 3 - moist
$ws2.Cells.Item(3,4).Value = "moist"
$ws2.Cells.Item(3,5).Value = 0.7

# Row 4 - Warm Temperate / dry
$ws2.Cells.Item(4,3).Value = "Warm Temperate"
$ws2.Cells.Item(4,4).Value = "-"
$ws2.Cells.Item(4,5).Value = 0.76

# Row 5 - moist
$ws2.Cells.Item(5,4).Value = "moist"
$ws2.Cells.Item(5,5).Value = 0.69

# Row 6 - Tropical / dry
$ws2.Cells.Item(6,3).Value = "Tropical"
$ws2.Cells.Item(6,4).Value = "-"
$ws2.Cells.Item(6,5).Value = 0.92

# Row 7 - moist/wet
$ws2.Cells.Item(7,4).Value = "moist/wet"
$ws2.Cells.Item(7,5).Value = 0.83

# Row 8 - rice
$ws2.Cells.Item(8,2).Value = "rice"
$ws2.Cells.Item(8,4).Value = 1.35

# Row 9 - perennial / temperate/boreal / mean lu factor
$ws2.Cells.Item(9,2).Value = "perennial"
$ws2.Cells.Item(9,3).Value = "temperate/boreal"
$ws2.Cells.Item(9,5).Value = 0.72
$ws2.Cells.Item(9,6).Value = "mean lu factor"

# Row 10 - tropical, average formula
$ws2.Cells.Item(10,3).Value = "tropical"
$ws2.Cells.Item(10,5).Value = 1.01
$ws2.Cells.Item(10,6).Formula = "=AVERAGE(E2:E10)"

# Row 12 - FMG / Reduced Till / Cool Temperate/Boreal / dry
$ws2.Cells.Item(12,1).Value = "FMG"
$ws2.Cells.Item(12,2).Value = "Reduced Till"
$ws2.Cells.Item(12,3).Value = "Cool Temperate/Boreal"
$ws2.Cells.Item(12,4).Value = "-"
$ws2.Cells.Item(12,5).Value = 0.98

# Row 13 - moist
$ws2.Cells.Item(13,4).Value = "moist"
$ws2.Cells.Item(13,5).Value = 1.04

# Row 14 - Warm Temperate / dry
$ws2.Cells.Item(14,3).Value = "Warm Temperate"
$ws2.Cells.Item(14,4).Value = "-"
$ws2.Cells.Item(14,5).Value = 0.99

# Row 15 - moist
$ws2.Cells.Item(15,4).Value = "moist"
$ws2.Cells.Item(15,5).Value = 1.05

# Row 16 - Tropical / dry
$ws2.Cells.Item(16,3).Value = "Tropical"
$ws2.Cells.Item(16,4).Value = "-"
$ws2.Cells.Item(16,5).Value = 0.99

# Row 17 - moist/wet, average formula
$ws2.Cells.Item(17,4).Value = "moist/wet"
$ws2.Cells.Item(17,5).Value = 1.04
$ws2.Cells.Item(17,6).Formula = "=AVERAGE(E12:E17)"

# Row 19 - No Till / Cool Temperate/Boreal / dry
$ws2.Cells.Item(19,2).Value = "No Till"
$ws2.Cells.Item(19,3).Value = "Cool Temperate/Boreal"
$ws2.Cells.Item(19,4).Value = "-"
$ws2.Cells.Item(19,5).Value = 1.03

# Row 20 - moist
$ws2.Cells.Item(20,4).Value = "moist"
$ws2.Cells.Item(20,5).Value = 1.09

# Row 21 - Warm Temperate / dry
$ws2.Cells.Item(21,3).Value = "Warm Temperate"
$ws2.Cells.Item(21,4).Value = "-"
$ws2.Cells.Item(21,5).Value = 1.04

# Row 22 - moist
$ws2.Cells.Item(22,4).Value = "moist"
$ws2.Cells.Item(22,5).Value = 1.1

# Row 23 - Tropical / dry
$ws2.Cells.Item(23,3).Value = "Tropical"
$ws2.Cells.Item(23,4).Value = "-"
$ws2.Cells.Item(23,5).Value = 1.04

# Row 24 - moist/wet, average formula
$ws2.Cells.Item(24,4).Value = "moist/wet"
$ws2.Cells.Item(24,5).Value = 1.1
$ws2.Cells.Item(24,6).Formula = "=AVERAGE(E19:E24)"

# Activate the new sheet and set its selection
$ws2.Activate()
$ws2.Range("E11").Select()
